$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (row 1, bold/border/date-format) into new column BB header cell
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)

# Copy left-column style (column A, bold/border/date-format) into new row 83 date cell
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)

# Set new column BB values (rows 1-83)
$ws.Range("BB1").Value = 45986
$ws.Range("BB2").Value = 0.8
$ws.Range("BB3").Value = -0.1
$ws.Range("BB4").Value = 0.4
$ws.Range("BB5").Value = -1.1
$ws.Range("BB6").Value = -1.7
$ws.Range("BB7").Value = -2.7
$ws.Range("BB8").Value = 1.7
$ws.Range("BB9").Value = -0.2
$ws.Range("BB10").Value = 1.3
$ws.Range("BB11").Value = -0.8
$ws.Range("BB12").Value = 0.8
$ws.Range("BB13").Value = 0.5
$ws.Range("BB14").Value = 0.2
$ws.Range("BB15").Value = 0.2
$ws.Range("BB16").Value = -0.7
$ws.Range("BB17").Value = 0.2000000000000001
$ws.Range("BB18").Value = 0.1
$ws.Range("BB19").Value = 0.7
$ws.Range("BB20").Value = 0.7
$ws.Range("BB21").Value = 0.5
$ws.Range("BB22").Value = -1.1
$ws.Range("BB23").Value = 0
$ws.Range("BB24").Value = 0
$ws.Range("BB25").Value = -0.6
$ws.Range("BB26").Value = 0.6
$ws.Range("BB27").Value = -0.1
$ws.Range("BB28").Value = -0.2
$ws.Range("BB29").Value = 0.5
$ws.Range("BB30").Value = -0.2
$ws.Range("BB31").Value = -0.1
$ws.Range("BB32").Value = 0.6
$ws.Range("BB33").Value = -0.5
$ws.Range("BB34").Value = -0.5
$ws.Range("BB35").Value = 0.1
$ws.Range("BB36").Value = 0.4
$ws.Range("BB37").Value = -0.4
$ws.Range("BB38").Value = -0.2
$ws.Range("BB39").Value = 0.5
$ws.Range("BB40").Value = -0.4
$ws.Range("BB41").Value = 0.4
$ws.Range("BB42").Value = 0.5
$ws.Range("BB43").Value = 0
$ws.Range("BB44").Value = -0.4
$ws.Range("BB45").Value = -0.7
$ws.Range("BB46").Value = -0.2
$ws.Range("BB47").Value = 0.2
$ws.Range("BB48").Value = -0.5
$ws.Range("BB49").Value = 0.8
$ws.Range("BB50").Value = -0.4
$ws.Range("BB51").Value = -1.1
$ws.Range("BB52").Value = -2.4
$ws.Range("BB53").Value = 3.8
$ws.Range("BB54").Value = 0.7
$ws.Range("BB55").Value = -0.9
$ws.Range("BB56").Value = -0.2
$ws.Range("BB57").Value = -0.1
$ws.Range("BB58").Value = 0.1
$ws.Range("BB59").Value = 0
$ws.Range("BB60").Value = -0.7
$ws.Range("BB61").Value = -0.4
$ws.Range("BB62").Value = 0.5
$ws.Range("BB63").Value = 0.8
$ws.Range("BB64").Value = -0.4
$ws.Range("BB65").Value = 0.4
$ws.Range("BB66").Value = 0.3
$ws.Range("BB67").Value = 0.1
$ws.Range("BB68").Value = -0.3
$ws.Range("BB69").Value = -1.1
$ws.Range("BB70").Value = -1.1
$ws.Range("BB71").Value = 0.4
$ws.Range("BB72").Value = -0.5
$ws.Range("BB73").Value = -0.3
$ws.Range("BB74").Value = -0.3
$ws.Range("BB75").Value = -0.3
$ws.Range("BB76").Value = -0.3
$ws.Range("BB77").Value = -0.3
$ws.Range("BB78").Value = -0.3
$ws.Range("BB79").Value = -0.3
$ws.Range("BB80").Value = -0.3
$ws.Range("BB81").Value = -0.3
$ws.Range("BB82").Value = -0.3
$ws.Range("BB83").Value = -0.3

# Set new row 83 date in column A
$ws.Range("A83").Value = 46934
